$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the alternating-shade row formatting down into the new row 29 ---
# Row 27 is the closest "odd" (shaded) template row; copy its formats down.
$ws.Range("A27:AK27").Copy() | Out-Null
$ws.Range("A29:AK29").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# The report-description cells on row 28 (P28 / AC28) pick up the wrap-text
# style used throughout the rest of the sheet once a new row follows them.
$ws.Range("P28").WrapText = $true
$ws.Range("AC28").WrapText = $true

# --- Populate the new row (new maintenance record, ticket #2025070735) ---
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "服務"
$ws.Range("C29").Value = 2025070735
$ws.Range("F29").Value = 4397
$ws.Range("G29").Value = "三重仁旺店"
$ws.Range("H29").Value = "新北市三重區"
$ws.Range("Q29").Value = "THILF04397"
$ws.Range("R29").Value = "新北一"
$ws.Range("S29").Value = "吳宗鴻"
$ws.Range("T29").Value = 1
$ws.Range("U29").Value = "已完工"
$ws.Range("V29").Value = "2025-07-04 14:46:07"
$ws.Range("W29").Value = "2025-07-04 13:50:00"
$ws.Range("X29").Value = "2025-07-04 14:45:00"
$ws.Range("Z29").Value = 0.9
$ws.Range("AB29").Value = "到場處理"
$ws.Range("AC29").Value = "PMQ3+STAR"
$ws.Range("AD29").Value = "O"
$ws.Range("AJ29").Value = "O"
$ws.Range("AK29").Value = "O"

# --- Grow the print area to cover the freshly-added row ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Report!Print_Area") {
        $n.RefersTo = "='Report'!`$A`$1:`$AK`$29"
    }
}

# --- Match the cursor position left behind by whoever added the row ---
$ws.Range("A29").Select() | Out-Null
